$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)
$jsonCol = $tbl.ListColumns.Item("JSON")

$newFormula = "=CONCATENATE(`"var s`",Table1[[#This Row],[id]],`" = {'publication':'`",Table1[[#This Row],[Newspaper]],`"','date': `",TEXT(Table1[[#This Row],[Formatted date]],`"yyyymmdd`"),`",'url':'`",Table1[[#This Row],[URL]],`"',debate:'`",Table1[[#This Row],[Debate]],`"'};`")"

$rng = $jsonCol.DataBodyRange
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    $rng.Cells.Item($i, 1).Formula = $newFormula
}
